$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 4.311421759901735
$ws.Range("C3").Value = 5.550740058157966
$ws.Range("E3").Value = 8.556235165581905
$ws.Range("C4").Value = 8.539761715194437
$ws.Range("E4").Value = 11.49866806009596
$ws.Range("C5").Value = 3.726709966233899
$ws.Range("E5").Value = 5.505010735461191
$ws.Range("C6").Value = -2.591890155624477
$ws.Range("C7").Value = 1.992279837313027
$ws.Range("E7").Value = 4.794307394577713
$ws.Range("C8").Value = 2.327315368300176
$ws.Range("C9").Value = -0.1842552248438545
$ws.Range("E9").Value = -1.471335939414509
$ws.Range("C10").Value = 2.026577416695763
$ws.Range("E10").Value = 1.772867810894829
$ws.Range("C11").Value = 2.854404831300794
$ws.Range("C12").Value = 1.119050958886225
$ws.Range("E12").Value = 1.764569308912711
$ws.Range("C13").Value = 2.543418408082077
$ws.Range("C14").Value = 1.807361459167756
$ws.Range("C15").Value = -0.004756022387275571
$ws.Range("C16").Value = 0.00555640981561023
$ws.Range("C17").Value = 0.9376318462105848
$ws.Range("E17").Value = 1.76760087614849
$ws.Range("C18").Value = 1.619232310145868
$ws.Range("E18").Value = 4.242290614103017
$ws.Range("C19").Value = 1.395050145291932
$ws.Range("E19").Value = 1.469426281897146
$ws.Range("C20").Value = 3.227100693237817
$ws.Range("C21").Value = 3.348613256881983
$ws.Range("C22").Value = -7.921833713986381
$ws.Range("C23").Value = 0.2827397234951956
$ws.Range("C24").Value = 4.451761251541475
$ws.Range("C25").Value = 1.493655572990393
$ws.Range("E25").Value = 2.056299233127357
$ws.Range("C26").Value = -0.2005234417569279
$ws.Range("C27").Value = 2.233381469093354
$ws.Range("C28").Value = 1.312922983354992
$ws.Range("E28").Value = 1.093658812337606
$ws.Range("C29").Value = 1.785401661837871
$ws.Range("C30").Value = 1.875259646256233
$ws.Range("C31").Value = 2.159588720360284
$ws.Range("E31").Value = 1.194436089410567
$ws.Range("C32").Value = 0.8260652760268661
$ws.Range("C33").Value = -3.195510012625546
$ws.Range("E33").Value = -8.595406358054735
$ws.Range("C34").Value = -1.620339334086651
$ws.Range("E34").Value = -7.034199224108983
$ws.Range("C35").Value = 1.905757520223461
$ws.Range("E35").Value = 0.9449497347025604
$ws.Range("C36").Value = -0.492451289919571
$ws.Range("E36").Value = -1.334179412476988
$ws.Range("C37").Value = -0.0458093613100008
$ws.Range("C38").Value = 0.4077310087939434
$ws.Range("E38").Value = 1.656917693245785
